$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.192.90'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.276.46'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.17'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.11'
$ws.Range("E6").Value = '  +3.62%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +4.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.418'
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.842.07'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.65'
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.189.33'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.270.00'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.63'
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '382.73'
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.77'
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  +1.28%  '
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("E26").Value = '  +6.89%  '
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("E29").Value = '  +2.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.99'
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.21'
$ws.Range("E32").Value = '  +5.83%  '
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.41'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.840'
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.70'
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.76'
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.62'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.64'
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.64'
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0692'
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.26'
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.632.28'
$ws.Range("E46").Value = '  -4.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '341.82'
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.20'
$ws.Range("E49").Value = '  +4.31%  '
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("E51").Value = '  -0.07%  '
